$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "HU_FertilityByYear" worksheet after "IT_FertilityByYear"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$huSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$huSheet.Name = "HU_FertilityByYear"

# Row 1: header "Fertility" then years 2011..2045 across columns B..AJ
$huSheet.Cells.Item(1, 1).Value = "Fertility"

$years = 2011..2045
for ($i = 0; $i -lt $years.Length; $i++) {
    $huSheet.Cells.Item(1, $i + 2).Value = $years[$i]
}

# Row 2: header "Value" then the fertility values across columns B..AJ
$huSheet.Cells.Item(2, 1).Value = "Value"

$values = @(
    41,
    44.666666666666664,
    45.000000000000007,
    48,
    48.333333333333336,
    51.000000000000007,
    51.333333333333336,
    51.666666666666664,
    51.666666666666664,
    53.000000000000007,
    53.666666666666671,
    52.000000000000007,
    52.000000000000007,
    52.000000000000007,
    54.740666666666669,
    54.740666666666669,
    54.740666666666669,
    54.740666666666669,
    54.740666666666669,
    55.614000000000004,
    55.614000000000004,
    55.614000000000004,
    55.614000000000004,
    55.614000000000004,
    56.212000000000003,
    56.212000000000003,
    56.212000000000003,
    56.212000000000003,
    56.212000000000003,
    56.612000000000002,
    56.612000000000002,
    56.612000000000002,
    56.612000000000002,
    56.612000000000002,
    56.87166666666667
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $huSheet.Cells.Item(2, $i + 2).Value = $values[$i]
}

# This new sheet becomes the active / selected tab.
$huSheet.Activate()

# ---------------------------------------------------------------------------
# 2. UK_FertilityByYear gains a selection (B1:BT1, active cell BT1)
# ---------------------------------------------------------------------------
$ukSheet = $wb.Worksheets.Item("UK_FertilityByYear")
$ukSheet.Activate()
$ukSheet.Range("B1:BT1").Select()

# ---------------------------------------------------------------------------
# 3. Re-activate the new HU sheet so it ends up as the selected/active tab
# ---------------------------------------------------------------------------
$huSheet.Activate()
